$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transactions")

# Rename portfolio value used in the "portfolio" column (I2:I7)
# from "FuturesPortfolioForPnLCalc" to "FuturesPortWithDiffCostBasis"
# to better reflect the example (SE-1843).
$range = $ws.Range("I2:I7")
foreach ($cell in $range.Cells) {
    $v = $cell.Value2
    if ($v -eq "FuturesPortfolioForPnLCalc") {
        $cell.Value = "FuturesPortWithDiffCostBasis"
    }
}

# Leave the cursor where the author left it when saving.
[void]$ws.Range("I9").Select()
